$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.851.15'
$ws.Range("E2").Value = '  -1.88%  '

# Row 3
$ws.Range("D3").Value = '1.802.40'
$ws.Range("E3").Value = '  -1.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = "'309.16"
$ws.Range("E5").Value = '  -1.78%  '

# Row 6
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("D7").Value = "'0.4640"
$ws.Range("E7").Value = '  +3.68%  '

# Row 8
$ws.Range("D8").Value = "'0.3694"
$ws.Range("E8").Value = '  -2.43%  '

# Row 9
$ws.Range("D9").Value = "'0.07363"
$ws.Range("E9").Value = '  -1.60%  '

# Row 10
$ws.Range("D10").Value = "'0.8681"
$ws.Range("E10").Value = '  -2.16%  '

# Row 11
$ws.Range("D11").Value = "'20.36"
$ws.Range("E11").Value = '  -3.43%  '

# Row 12
$ws.Range("D12").Value = '1.778.79'
$ws.Range("E12").Value = '  -2.67%  '

# Row 13
$ws.Range("D13").Value = "'5.354"
$ws.Range("E13").Value = '  -1.86%  '

# Row 14
$ws.Range("D14").Value = "'92.14"
$ws.Range("E14").Value = '  -1.85%  '

# Row 15
$ws.Range("E15").Value = '  -3.61%  '

# Row 16
$ws.Range("D16").Value = "'0.07030"
$ws.Range("E16").Value = '  -1.26%  '

# Row 17
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  -0.02%  '

# Row 18
$ws.Range("D18").Value = "'0.000008706"
$ws.Range("E18").Value = '  -1.04%  '

# Row 19
$ws.Range("E19").Value = '  +0.06%  '

# Row 20
$ws.Range("E20").Value = '  -3.43%  '

# Row 21
$ws.Range("D21").Value = '26.850.35'
$ws.Range("E21").Value = '  -1.88%  '

# Row 22
$ws.Range("E22").Value = '  -2.07%  '

# Row 23
$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = '  -3.72%  '

# Row 24
$ws.Range("D24").Value = '2.078.91'
$ws.Range("E24").Value = '  +1.02%  '

# Row 25
$ws.Range("D25").Value = "'1.902"
$ws.Range("E25").Value = '  -3.13%  '

# Row 26
$ws.Range("D26").Value = "'151.33"
$ws.Range("E26").Value = '  -0.08%  '

# Row 27
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = '  -1.89%  '

# Row 28
$ws.Range("D28").Value = "'2.135"
$ws.Range("E28").Value = '  -7.50%  '

# Row 29
$ws.Range("E29").Value = '  -3.30%  '

# Row 30
$ws.Range("D30").Value = "'115.93"
$ws.Range("E30").Value = '  -1.79%  '

# Row 31
$ws.Range("D31").Value = "'0.08912"
$ws.Range("E31").Value = '  +0.24%  '

# Row 32
$ws.Range("D32").Value = "'0.7590"
$ws.Range("E32").Value = '  -4.08%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'1.150"
$ws.Range("E33").Value = '  -4.87%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = "'2.930"
$ws.Range("E34").Value = '  +0.23%  '

# Row 35
$ws.Range("D35").Value = "'4.454"
$ws.Range("E35").Value = '  -3.29%  '

# Row 36
$ws.Range("D36").Value = "'0.9996"
$ws.Range("E36").Value = '  -0.05%  '

# Row 37
$ws.Range("E37").Value = '  -0.72%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.01948"
$ws.Range("E38").Value = '  -2.24%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'0.05245"
$ws.Range("E39").Value = '  -1.34%  '

# Row 40
$ws.Range("E40").Value = '  +2.07%  '

# Row 41
$ws.Range("D41").Value = "'7.201"
$ws.Range("E41").Value = '  -1.79%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'2.364"
$ws.Range("E42").Value = '  +1.67%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = "'0.5285"
$ws.Range("E43").Value = '  -1.31%  '

# Row 44
$ws.Range("E44").Value = '  -3.53%  '

# Row 45
$ws.Range("D45").Value = "'8.487"
$ws.Range("E45").Value = '  -2.31%  '

# Row 46
$ws.Range("D46").Value = "'0.5001"
$ws.Range("E46").Value = '  -2.36%  '

# Row 47
$ws.Range("D47").Value = "'10.27"
$ws.Range("E47").Value = '  -3.95%  '

# Row 48
$ws.Range("E48").Value = '  -1.30%  '

# Row 49
$ws.Range("D49").Value = "'0.9995"
$ws.Range("E49").Value = '  -0.03%  '

# Row 50
$ws.Range("D50").Value = "'1.662"

# Row 51
$ws.Range("D51").Value = "'0.06285"
$ws.Range("E51").Value = '  -1.93%  '
